# Workbook was edited in real Excel/openpyxl: the "Child" table addresses
# were regenerated, and the active sheet/selection moved from "Child"
# (J21) to the newly built "School" object (K14), with the School sheet
# view picking up the same 145% zoom the Child sheet already used.

$wb = $excel.ActiveWorkbook

# --- Update the randomly generated "address" column on the Child sheet ---
$child = $wb.Worksheets.Item("Child")
$child.Range("D2").Value  = "9,-2"
$child.Range("D3").Value  = "2,-1"
$child.Range("D4").Value  = "-7,-3"
$child.Range("D5").Value  = "8,6"
$child.Range("D6").Value  = "2,0"
$child.Range("D7").Value  = "-3,5"
$child.Range("D8").Value  = "2,-1"
$child.Range("D9").Value  = "8,3"
$child.Range("D10").Value = "-5,-10"
$child.Range("D11").Value = "-9,-6"
$child.Range("D12").Value = "7,6"
$child.Range("D13").Value = "8,0"
$child.Range("D14").Value = "3,5"
$child.Range("D15").Value = "5,-7"
$child.Range("D16").Value = "-4,9"
$child.Range("D17").Value = "-7,-6"
$child.Range("D18").Value = "0,-9"
$child.Range("D19").Value = "-3,-8"
$child.Range("D20").Value = "-3,2"
$child.Range("D21").Value = "1,-2"

# Keep the Child sheet's own selection where it was (J21) -- only the
# active/tab-selected sheet is moving, not this sheet's cursor.
$child.Range("J21").Select() | Out-Null

# --- Move the active window to the School sheet ---
$school = $wb.Worksheets.Item("School")
$school.Activate() | Out-Null
$school.Range("K14").Select() | Out-Null
$excel.ActiveWindow.Zoom = 145
